$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the existing row 307.
# This automatically shifts the existing rows 307..381 down to 310..384,
# matching the rest of the diff (which is just the old rows re-numbered).
$ws.Rows.Item(307).Resize(3).Insert()

# Common/constant values shared by every data row in this sheet.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100112023
$categoria = "Brócoli"
$variedad  = "Sin especificar"
$calidad   = "Primera"
$unidad    = "`$/unidad"
$kgUnid    = 1
$clasif    = "Hortaliza"

function Set-DataRow {
    param(
        [int]$RowNum,
        [int]$Fecha,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Origen,
        [double]$PrecioKg
    )

    $ws.Cells.Item($RowNum, 1).Value  = $mercadoId
    $ws.Cells.Item($RowNum, 2).Value  = $mercado
    $ws.Cells.Item($RowNum, 3).Value  = $region
    $ws.Cells.Item($RowNum, 4).Value  = $Fecha
    $ws.Cells.Item($RowNum, 5).Value  = $codreg
    $ws.Cells.Item($RowNum, 6).Value  = $catId
    $ws.Cells.Item($RowNum, 7).Value  = $categoria
    $ws.Cells.Item($RowNum, 8).Value  = $variedad
    $ws.Cells.Item($RowNum, 9).Value  = $calidad
    $ws.Cells.Item($RowNum, 10).Value = $Volumen
    $ws.Cells.Item($RowNum, 11).Value = $PrecioMin
    $ws.Cells.Item($RowNum, 12).Value = $PrecioMax
    $ws.Cells.Item($RowNum, 13).Value = $PrecioProm
    $ws.Cells.Item($RowNum, 14).Value = $unidad
    $ws.Cells.Item($RowNum, 15).Value = $Origen
    $ws.Cells.Item($RowNum, 16).Value = $PrecioKg
    $ws.Cells.Item($RowNum, 17).Value = $kgUnid
    $ws.Cells.Item($RowNum, 18).Value = $clasif
}

# New row 307
Set-DataRow 307 44642 310 1400 1400 1400 "Provincia de Cautín" 1400

# New row 308
Set-DataRow 308 44642 850 1000 1100 1053 "Región Metropolitana" 1053

# New row 309
Set-DataRow 309 44642 750 1200 1200 1200 "Región del Maule" 1200
